$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Rousset2004): tntRuns & rRuns complete (numeric 1 instead of placeholder "zzz"),
# and pdfOut now also marked complete.
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 1

# Row 8: rRuns placeholder "zzz" added.
$ws.Range("J8").Value = "zzz"

# Row 13: charType marked complete.
$ws.Range("H13").Value = 1
